# Update gh-pages to output generated at 456a3b4
# Applies the F-column ("想去人数" / want-to-go count) updates across the
# four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F13").Value = 2408
$ws1.Range("F18").Value = 529
$ws1.Range("F21").Value = 9982
$ws1.Range("F22").Value = 5950
$ws1.Range("F26").Value = 138
$ws1.Range("F28").Value = 3529
$ws1.Range("F31").Value = 457
$ws1.Range("F34").Value = 236
$ws1.Range("F35").Value = 224
$ws1.Range("F36").Value = 4811
$ws1.Range("F37").Value = 18
$ws1.Range("F39").Value = 146
$ws1.Range("F40").Value = 22
$ws1.Range("F41").Value = 60

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 3525

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1557

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1557
$ws4.Range("F13").Value = 2408
$ws4.Range("F22").Value = 529
$ws4.Range("F25").Value = 9982
$ws4.Range("F26").Value = 3525
$ws4.Range("F31").Value = 138
$ws4.Range("F33").Value = 3529
$ws4.Range("F36").Value = 457
$ws4.Range("F38").Value = 236
$ws4.Range("F40").Value = 224
$ws4.Range("F41").Value = 4811
$ws4.Range("F44").Value = 60
